$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("second")

# Capture all current row values (using Value2, which works reliably in this runtime)
$a1 = $ws.Range("A1").Value2
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2

$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2

# Shift rows 2-5 up into rows 1-4, and move the original row 1 into row 5
$ws.Range("A1").Value = $a2
$ws.Range("B1").Value = $b2
$ws.Range("C1").Value = $c2

$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3

$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("C3").Value = $c4

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("C4").Value = $c5

$ws.Range("A5").Value = $a1
$ws.Range("B5").Value = $b1
$ws.Range("C5").Value = $c1

# Update the selection to match the target state
$ws.Range("A8").Select()

$wb.Save()
